# Update "Estimates and Actuals.xlsx":
#  - Weekly sheet: add Week-4 ("F" column) actuals for a few tasks, and add a
#    new task row "Loc Implementation/refinement" with a Week-4 actual.
#  - Summary sheet: some Estimate/Actual cells on the summary tab are edited
#    directly (F3 and F5 become typed-in values instead of formulas pulling
#    from the Weekly sheet), plus a couple of other Actual entries.

$wb = $excel.ActiveWorkbook

$wsWeekly  = $wb.Worksheets.Item("Weekly")
$wsSummary = $wb.Worksheets.Item("Summary")

# --- Weekly sheet updates ---
$wsWeekly.Range("F6").Value = 0.75
$wsWeekly.Range("F8").Value = 6
$wsWeekly.Range("F10").Value = 4

# New row 11: a new task with its actual hours
$wsWeekly.Range("B11").Value = "Loc Implementation/refinement"
$wsWeekly.Range("F11").Value = 6

# --- Summary sheet updates ---
# F3 and F5 were formulas referencing the Weekly sheet; they are now
# overwritten with plain typed-in values.
$wsSummary.Range("F3").Value = 3.75
$wsSummary.Range("F5").Value = 14

$wsSummary.Range("F9").Value = 2
$wsSummary.Range("F14").Value = 4
$wsSummary.Range("F17").Value = 16

# --- Match final cursor/selection positions ---
$wsWeekly.Range("F10").Select()
$wsSummary.Select()
$wsSummary.Range("H18").Select()
